# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = "42.937.94"
$ws.Range("E2").Value = "  +0.85%  "

# Row 3
$ws.Range("D3").Value = "2.295.67"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'316.33"
$ws.Range("E5").Value = "  +0.08%  "

# Row 6
$ws.Range("D6").Value = "'104.45"
$ws.Range("E6").Value = "  +0.23%  "

# Row 7
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  -1.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  -1.56%  "

# Row 10
$ws.Range("D10").Value = "'39.49"
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
$ws.Range("E11").Value = "  -0.65%  "

# Row 12
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "  +1.50%  "

# Row 13
$ws.Range("E13").Value = "  +2.19%  "

# Row 14
$ws.Range("E14").Value = "  +4.38%  "

# Row 15
$ws.Range("D15").Value = "'15.34"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").Value = "2.642.36"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17
$ws.Range("D17").Value = "2.309.33"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").Value = "42.814.53"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").Value = "'13.90"
$ws.Range("E20").Value = "  +25.48%  "

# Row 21
$ws.Range("E21").Value = "  -0.56%  "

# Row 22
$ws.Range("D22").Value = "'73.99"
$ws.Range("E22").Value = "  +0.98%  "

# Row 23
$ws.Range("D23").Value = "'3.56"
$ws.Range("E23").Value = "  +0.85%  "

# Row 24
$ws.Range("E24").Value = "  -4.94%  "

# Row 25
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  -3.05%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.18%  "

# Row 27
$ws.Range("D27").Value = "'10.92"
$ws.Range("E27").Value = "  +0.85%  "

# Row 28
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'7.08"
$ws.Range("E28").Value = "  +20.04%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "  -2.76%  "

# Row 30
$ws.Range("D30").Value = "'22.38"
$ws.Range("E30").Value = "  -1.72%  "

# Row 31
$ws.Range("D31").Value = "'37.61"
$ws.Range("E31").Value = "  +5.36%  "

# Row 32
$ws.Range("D32").Value = "'166.64"
$ws.Range("E32").Value = "  +1.29%  "

# Row 33
$ws.Range("D33").Value = "'0.0875"
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$ws.Range("E34").Value = "  -4.08%  "

# Row 35
$ws.Range("E35").Value = "  -0.63%  "

# Row 36
$ws.Range("E36").Value = "  -1.37%  "

# Row 37
$ws.Range("D37").Value = "'4.57"
$ws.Range("E37").Value = "  -0.33%  "

# Row 38
$ws.Range("D38").Value = "'0.0351"
$ws.Range("E38").Value = "  -5.03%  "

# Row 39
$ws.Range("D39").Value = "'3.85"
$ws.Range("E39").Value = "  +3.17%  "

# Row 40
$ws.Range("E40").Value = "  -2.59%  "

# Row 41
$ws.Range("E41").Value = "  +5.41%  "

# Row 42
$ws.Range("E42").Value = "  +1.21%  "

# Row 43
$ws.Range("D43").Value = "'69.71"
$ws.Range("E43").Value = "  -0.34%  "

# Row 44
$ws.Range("E44").Value = "  +0.21%  "

# Row 45
$ws.Range("D45").Value = "'92.91"
$ws.Range("E45").Value = "  -1.37%  "

# Row 46
$ws.Range("D46").Value = "'12.22"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("D47").Value = "'114.53"
$ws.Range("E47").Value = "  +1.37%  "

# Row 48
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'80.41"
$ws.Range("E48").Value = "  -2.64%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.724.35"
$ws.Range("E49").Value = "  +8.41%  "

# Row 50
$ws.Range("D50").Value = "'8.78"
$ws.Range("E50").Value = "  -1.37%  "

# Row 51
$ws.Range("D51").Value = "'5.16"
$ws.Range("E51").Value = "  +1.89%  "
